# Sync automático del tracker - add a "Results" sheet after "Predictions"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Clone the Predictions sheet so the new sheet inherits identical styles
# (header style, borders, etc.) instead of Excel synthesizing brand-new
# style/font entries for freshly-typed formatting.
$ws1.Copy($null, $ws1)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Results"

# Drop the cloned sample rows and the extra (Predictions-only) columns,
# leaving just a single header row A1:J1.
$newSheet.Rows("2:5").Delete()
$newSheet.Columns("K:P").Delete()

# Overwrite the header labels for the Results tracker.
$headers = @("Date", "Liga", "Local", "Visitante", "Resultado_Real", "Predicción", "Acierto", "Profit", "ROI", "Fecha_Partido")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Keep the original active sheet selected, as before the edit.
$ws1.Select()
